# Adds a new "Mae pc от old " column (new column I) to each of the three
# report sheets, renames the existing Mae column header to "Mae old от pc",
# pushes the old "Тип данных" column to J, fills the new I column with a
# copy of the H (Mae) values, and reverses the data-row order on the
# "Negative Correlation" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Positive Correlation", "Negative Correlation", "General")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the last used row/column on this sheet.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column

    # Insert a new blank column at I; the existing I ("Тип данных") and
    # everything after it shifts right to become J, K, ...
    $ws.Columns.Item(9).Insert()

    # Header row: H1 is renamed, new I1 gets the new header text.
    $ws.Cells.Item(1, 8).Value2 = "Mae old от pc"
    $ws.Cells.Item(1, 9).Value2 = "Mae pc от old "

    # Carry over the formatting used by the other header cells (bold/border/
    # centered) onto the freshly inserted I1 header cell.
    $headerSrc = $ws.Cells.Item(1, 8)
    $headerDst = $ws.Cells.Item(1, 9)
    $headerDst.Font.Bold = $headerSrc.Font.Bold
    $headerDst.HorizontalAlignment = $headerSrc.HorizontalAlignment
    $headerDst.VerticalAlignment = $headerSrc.VerticalAlignment
    $headerDst.Borders.LineStyle = $headerSrc.Borders.LineStyle

    # Fill the new I column (data rows) with a copy of the H (Mae) values.
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 9).Value2 = $ws.Cells.Item($r, 8).Value2
    }

    if ($sheetName -eq "Negative Correlation") {
        # Reverse the order of the data rows (row 2..lastRow), now that the
        # sheet has columns A..J.
        $newLastCol = $lastCol + 1
        $fullRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, $newLastCol))
        $vals = $fullRange.Value2

        $n = $vals.GetLength(0)
        $m = $vals.GetLength(1)
        $rev = New-Object 'object[,]' $n, $m
        for ($i = 1; $i -le $n; $i++) {
            for ($j = 1; $j -le $m; $j++) {
                $rev[$i - 1, $j - 1] = $vals[$n - $i + 1, $j]
            }
        }
        $fullRange.Value2 = $rev
    }
}
